# Daily update: decrement the "剩余" (remaining) days counter for every
# shop row. When a row's remaining count would drop to 0 (i.e. it was 1),
# the cycle restarts: remaining resets to the row's total day count (D)
# and the start date (F) advances by 7 days. Rows whose start date is not
# a valid yyyyMMdd date (data anomalies) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 99 }

for ($r = 2; $r -le $lastRow; $r++) {
    $total = $ws.Cells.Item($r, 4).Value2
    $remaining = $ws.Cells.Item($r, 5).Value2
    $startRaw = $ws.Cells.Item($r, 6).Value2

    if ($remaining -eq $null -or $total -eq $null) { continue }

    try {
        $startDate = [datetime]::ParseExact([string]$startRaw, "yyyyMMdd", $null)
    } catch {
        continue
    }

    if ($remaining -eq 1) {
        $ws.Cells.Item($r, 5).Value = $total
        $newStart = $startDate.AddDays(7)
        $ws.Cells.Item($r, 6).Value = [int]$newStart.ToString("yyyyMMdd")
    } else {
        $ws.Cells.Item($r, 5).Value = $remaining - 1
    }
}
